$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Row 28: relabel the existing vimrc snapshot entry with a date ---
$ws.Range("B28").Value = "snapshot/ 20170526"

# --- Row 32 (new): Tagbar / Add markdown support ---
$row32C = @'
# Add markdown support for tagbar
1. download https://github.com/jszakmeister/markdown2ctags to somewhere.
1. add below snippet to .vimrc:
" Add support for markdown files in tagbar.
let g:tagbar_type_markdown = {
    \ 'ctagstype': 'markdown',
    \ 'ctagsbin' : '/path/to/markdown2ctags.py',
    \ 'ctagsargs' : '-f - --sort=yes',
    \ 'kinds' : [
        \ 's:sections',
        \ 'i:images'
    \ ],
    \ 'sro' : '|',
    \ 'kind2scope' : {
        \ 's' : 'section',
    \ },
    \ 'sort': 0,
\ }
'@
$ws.Range("A32").Value = "Tagbar"
$ws.Range("B32").Value = "Add markdown support"
$ws.Range("C32").Value = $row32C
$ws.Rows.Item(32).RowHeight = 95.25

# --- Row 33 (new): vimrc / snapshot/ 20170526 (updated content with markdown support) ---
$row33C = @'
" myles put 2017-05-09 on indent tab to space
filetype plugin indent on
set tabstop=4                   "show existing tab with 4 spaces width
set shiftwidth=4                " when indenting with '>', use 4 spaces width
set expandtab                   " On pressing tab, insert 4 spaces
" myles put 2017-05-09 on set default colorscheme
colorscheme elflord
" Install pathogen 2017-05-23
execute pathogen#infect()
syntax on
" 2017-05-23 set option
"set nowrap
" 2017-05-24 set option
set incsearch
set hlsearch
" 2017-05-24 Open Nerdtree when vim open with no input
"autocmd StdinReadPre * let s:std_in=1
"autocmd VimEnter * if argc() == 0 && !exists("s:std_in") | NERDTree | endif
" 2017-05-24 Create shortcut for NERDTree
map <C-n> :NERDTreeToggle<CR>
" 2017-05-25 Change the directory arrow
let g:NERDTreeDirArrowExpandable="+"
let g:NERDTreeDirArrowCollapsible="-"
" 2017-05-25 Change the directory arrow
autocmd BufWinLeave *.* mkview
"autocmd BufWinEnter *.* silent loadview
" 2017-05-25 Syntastic plugin
" autocmd BufWinEnter *.* silent SyntasticToggleMode
" set statusline+=%#warningmsg#
" set statusline+=%{SyntasticStatuslineFlag()}
" set statusline+=%*
" let g:syntastic_always_populate_loc_list = 1
" let g:syntastic_auto_loc_list = 1
" let g:syntastic_check_on_open = 1
" let g:syntastic_check_on_wq = 0
" 2017-05-26 Syntax folding autostart
set foldmethod=syntax
set foldlevel=1
" 2017-06-13 Shortcut for Tagbar/ learnt that <ENTER> also bind to map
map <C-m> :TagbarToggle<CR>
" 2017-06-26 Add support for markdown files in tagbar.
let g:tagbar_type_markdown = {
    \ 'ctagstype': 'markdown',
    \ 'ctagsbin' : '~/.vim/support-script/markdown2ctags/markdown2ctags.py',
    \ 'ctagsargs' : '-f - --sort=yes',
    \ 'kinds' : [
        \ 's:sections',
        \ 'i:images'
    \ ],
    \ 'sro' : '|',
    \ 'kind2scope' : {
        \ 's' : 'section',
    \ },
    \ 'sort': 0,
\ }

'@
$ws.Range("A33").Value = "vimrc"
$ws.Range("B33").Value = "snapshot/ 20170526"
$ws.Range("C33").Value = $row33C
$ws.Rows.Item(33).RowHeight = 101.25

